# Auto-generated: update BRVM Recommandations + Top_YTD data (GitHub Actions MAJ automatique)
$wb = $excel.ActiveWorkbook
$wsReco = $wb.Worksheets.Item("Recommandations")
$wsYtd = $wb.Worksheets.Item("Top_YTD")

# --- Recommandations sheet: rows 2-47 (A:title, B:hausse, C:baisse, D:variation totale, E:derniere variation, F:reco, G:strategie) ---
$wsReco.Cells.Item(2,1).Value = "BRVM - SERVICES PUBLICS"
$wsReco.Cells.Item(2,2).Value = 0
$wsReco.Cells.Item(2,3).Value = 8
$wsReco.Cells.Item(2,4).Value = 3440.49
$wsReco.Cells.Item(2,5).Value = 113.66
$wsReco.Cells.Item(2,6).Value = "🟡 Observer"
$wsReco.Cells.Item(2,7).Value = "➖ Neutre"

$wsReco.Cells.Item(3,1).Value = "AIR LIQUIDE CI"
$wsReco.Cells.Item(3,2).Value = 0
$wsReco.Cells.Item(3,3).Value = 4
$wsReco.Cells.Item(3,4).Value = 2805
$wsReco.Cells.Item(3,5).Value = 710
$wsReco.Cells.Item(3,6).Value = "🟡 Observer"
$wsReco.Cells.Item(3,7).Value = "➖ Neutre"

$wsReco.Cells.Item(4,1).Value = "NEI-CEDA CI"
$wsReco.Cells.Item(4,2).Value = 0
$wsReco.Cells.Item(4,3).Value = 4
$wsReco.Cells.Item(4,4).Value = 2740
$wsReco.Cells.Item(4,5).Value = 685
$wsReco.Cells.Item(4,6).Value = "🟡 Observer"
$wsReco.Cells.Item(4,7).Value = "➖ Neutre"

$wsReco.Cells.Item(5,1).Value = "BRVM - AUTRES SECTEURS"
$wsReco.Cells.Item(5,2).Value = 0
$wsReco.Cells.Item(5,3).Value = 4
$wsReco.Cells.Item(5,4).Value = 2448.08
$wsReco.Cells.Item(5,5).Value = 610.03
$wsReco.Cells.Item(5,6).Value = "🟡 Observer"
$wsReco.Cells.Item(5,7).Value = "➖ Neutre"

$wsReco.Cells.Item(6,1).Value = "BRVM - DISTRIBUTION"
$wsReco.Cells.Item(6,2).Value = 0
$wsReco.Cells.Item(6,3).Value = 4
$wsReco.Cells.Item(6,4).Value = 2169.93
$wsReco.Cells.Item(6,5).Value = 564.45
$wsReco.Cells.Item(6,6).Value = "🟡 Observer"
$wsReco.Cells.Item(6,7).Value = "➖ Neutre"

$wsReco.Cells.Item(7,1).Value = "BRVM - AGRICULTURE"
$wsReco.Cells.Item(7,2).Value = 0
$wsReco.Cells.Item(7,3).Value = 4
$wsReco.Cells.Item(7,4).Value = 1471.26
$wsReco.Cells.Item(7,5).Value = 370.6
$wsReco.Cells.Item(7,6).Value = "🟡 Observer"
$wsReco.Cells.Item(7,7).Value = "➖ Neutre"

$wsReco.Cells.Item(8,1).Value = "BRVM - TRANSPORT"
$wsReco.Cells.Item(8,2).Value = 0
$wsReco.Cells.Item(8,3).Value = 4
$wsReco.Cells.Item(8,4).Value = 1429.47
$wsReco.Cells.Item(8,5).Value = 354.92
$wsReco.Cells.Item(8,6).Value = "🟡 Observer"
$wsReco.Cells.Item(8,7).Value = "➖ Neutre"

$wsReco.Cells.Item(9,1).Value = "BRVM - CONSOMMATION DISCRETIONNAIRE"
$wsReco.Cells.Item(9,2).Value = 0
$wsReco.Cells.Item(9,3).Value = 4
$wsReco.Cells.Item(9,4).Value = 796.81
$wsReco.Cells.Item(9,5).Value = 210.16
$wsReco.Cells.Item(9,6).Value = "🟡 Observer"
$wsReco.Cells.Item(9,7).Value = "➖ Neutre"

$wsReco.Cells.Item(10,1).Value = "BRVM - FINANCES"
$wsReco.Cells.Item(10,2).Value = 0
$wsReco.Cells.Item(10,3).Value = 4
$wsReco.Cells.Item(10,4).Value = 568.75
$wsReco.Cells.Item(10,5).Value = 143.73
$wsReco.Cells.Item(10,6).Value = "🟡 Observer"
$wsReco.Cells.Item(10,7).Value = "➖ Neutre"

$wsReco.Cells.Item(11,1).Value = "BRVM-PRESTIGE"
$wsReco.Cells.Item(11,2).Value = 0
$wsReco.Cells.Item(11,3).Value = 4
$wsReco.Cells.Item(11,4).Value = 567.84
$wsReco.Cells.Item(11,5).Value = 143.63
$wsReco.Cells.Item(11,6).Value = "🟡 Observer"
$wsReco.Cells.Item(11,7).Value = "➖ Neutre"

$wsReco.Cells.Item(12,1).Value = "BRVM - SERVICES FINANCIERS"
$wsReco.Cells.Item(12,2).Value = 0
$wsReco.Cells.Item(12,3).Value = 4
$wsReco.Cells.Item(12,4).Value = 558.97
$wsReco.Cells.Item(12,5).Value = 141.26
$wsReco.Cells.Item(12,6).Value = "🟡 Observer"
$wsReco.Cells.Item(12,7).Value = "➖ Neutre"

$wsReco.Cells.Item(13,1).Value = "BRVM - INDUSTRIELS"
$wsReco.Cells.Item(13,2).Value = 0
$wsReco.Cells.Item(13,3).Value = 4
$wsReco.Cells.Item(13,4).Value = 509.44
$wsReco.Cells.Item(13,5).Value = 127.69
$wsReco.Cells.Item(13,6).Value = "🟡 Observer"
$wsReco.Cells.Item(13,7).Value = "➖ Neutre"

$wsReco.Cells.Item(14,1).Value = "BRVM - ENERGIE"
$wsReco.Cells.Item(14,2).Value = 0
$wsReco.Cells.Item(14,3).Value = 4
$wsReco.Cells.Item(14,4).Value = 433.49
$wsReco.Cells.Item(14,5).Value = 107.92
$wsReco.Cells.Item(14,6).Value = "🟡 Observer"
$wsReco.Cells.Item(14,7).Value = "➖ Neutre"

$wsReco.Cells.Item(15,1).Value = "BRVM - TELECOMMUNICATIONS"
$wsReco.Cells.Item(15,2).Value = 0
$wsReco.Cells.Item(15,3).Value = 4
$wsReco.Cells.Item(15,4).Value = 389.31
$wsReco.Cells.Item(15,5).Value = 97.27
$wsReco.Cells.Item(15,6).Value = "🟡 Observer"
$wsReco.Cells.Item(15,7).Value = "➖ Neutre"

$wsReco.Cells.Item(16,1).Value = "BRVM - CONSOMMATION DE BASE         (**)"
$wsReco.Cells.Item(16,2).Value = 0
$wsReco.Cells.Item(16,3).Value = 2
$wsReco.Cells.Item(16,4).Value = 380.28
$wsReco.Cells.Item(16,5).Value = 190.31
$wsReco.Cells.Item(16,6).Value = "🟡 Observer"
$wsReco.Cells.Item(16,7).Value = "➖ Neutre"

$wsReco.Cells.Item(17,1).Value = "BRVM - INDUSTRIE                 (**)"
$wsReco.Cells.Item(17,2).Value = 0
$wsReco.Cells.Item(17,3).Value = 1
$wsReco.Cells.Item(17,4).Value = 215.1
$wsReco.Cells.Item(17,5).Value = 215.1
$wsReco.Cells.Item(17,6).Value = "🟡 Observer"
$wsReco.Cells.Item(17,7).Value = "➖ Neutre"

$wsReco.Cells.Item(18,1).Value = "BRVM - INDUSTRIE         (**)"
$wsReco.Cells.Item(18,2).Value = 0
$wsReco.Cells.Item(18,3).Value = 1
$wsReco.Cells.Item(18,4).Value = 214.39
$wsReco.Cells.Item(18,5).Value = 214.39
$wsReco.Cells.Item(18,6).Value = "🟡 Observer"
$wsReco.Cells.Item(18,7).Value = "➖ Neutre"

$wsReco.Cells.Item(19,1).Value = "BRVM-PRINCIPAL            (**)"
$wsReco.Cells.Item(19,2).Value = 0
$wsReco.Cells.Item(19,3).Value = 1
$wsReco.Cells.Item(19,4).Value = 205.89
$wsReco.Cells.Item(19,5).Value = 205.89
$wsReco.Cells.Item(19,6).Value = "🟡 Observer"
$wsReco.Cells.Item(19,7).Value = "➖ Neutre"

$wsReco.Cells.Item(20,1).Value = "BRVM-PRINCIPAL                   (**)"
$wsReco.Cells.Item(20,2).Value = 0
$wsReco.Cells.Item(20,3).Value = 1
$wsReco.Cells.Item(20,4).Value = 205.64
$wsReco.Cells.Item(20,5).Value = 205.64
$wsReco.Cells.Item(20,6).Value = "🟡 Observer"
$wsReco.Cells.Item(20,7).Value = "➖ Neutre"

$wsReco.Cells.Item(21,1).Value = "SAFCA CI (SAFC)"
$wsReco.Cells.Item(21,2).Value = 4
$wsReco.Cells.Item(21,3).Value = 0
$wsReco.Cells.Item(21,4).Value = 24.83
$wsReco.Cells.Item(21,5).Value = 7.43
$wsReco.Cells.Item(21,6).Value = "🟢 Achat"
$wsReco.Cells.Item(21,7).Value = "✅ Renforcer"

$wsReco.Cells.Item(22,1).Value = "CFAO MOTORS CI (CFAC)"
$wsReco.Cells.Item(22,2).Value = 3
$wsReco.Cells.Item(22,3).Value = 0
$wsReco.Cells.Item(22,4).Value = 19.88
$wsReco.Cells.Item(22,5).Value = 7.45
$wsReco.Cells.Item(22,6).Value = "🟢 Achat"
$wsReco.Cells.Item(22,7).Value = "✅ Renforcer"

$wsReco.Cells.Item(23,1).Value = "SERVAIR ABIDJAN CI (ABJC)"
$wsReco.Cells.Item(23,2).Value = 2
$wsReco.Cells.Item(23,3).Value = 1
$wsReco.Cells.Item(23,4).Value = 8
$wsReco.Cells.Item(23,5).Value = 7.5
$wsReco.Cells.Item(23,6).Value = "🟡 Observer"
$wsReco.Cells.Item(23,7).Value = "👀 À surveiller"

$wsReco.Cells.Item(24,1).Value = "SETAO CI (STAC)"
$wsReco.Cells.Item(24,2).Value = 1
$wsReco.Cells.Item(24,3).Value = 0
$wsReco.Cells.Item(24,4).Value = 7.5
$wsReco.Cells.Item(24,5).Value = 7.5
$wsReco.Cells.Item(24,6).Value = "🟡 Observer"
$wsReco.Cells.Item(24,7).Value = "➖ Neutre"

$wsReco.Cells.Item(25,1).Value = "SUCRIVOIRE (SCRC)"
$wsReco.Cells.Item(25,2).Value = 1
$wsReco.Cells.Item(25,3).Value = 0
$wsReco.Cells.Item(25,4).Value = 5.68
$wsReco.Cells.Item(25,5).Value = 5.68
$wsReco.Cells.Item(25,6).Value = "🟡 Observer"
$wsReco.Cells.Item(25,7).Value = "➖ Neutre"

$wsReco.Cells.Item(26,1).Value = "ECOBANK COTE D''IVOIRE (ECOC)"
$wsReco.Cells.Item(26,2).Value = 1
$wsReco.Cells.Item(26,3).Value = 0
$wsReco.Cells.Item(26,4).Value = 5.6
$wsReco.Cells.Item(26,5).Value = 5.6
$wsReco.Cells.Item(26,6).Value = "🟡 Observer"
$wsReco.Cells.Item(26,7).Value = "➖ Neutre"

$wsReco.Cells.Item(27,1).Value = "ECOBANK TRANS. INCORP. TG (ETIT)"
$wsReco.Cells.Item(27,2).Value = 1
$wsReco.Cells.Item(27,3).Value = 0
$wsReco.Cells.Item(27,4).Value = 4.76
$wsReco.Cells.Item(27,5).Value = 4.76
$wsReco.Cells.Item(27,6).Value = "🟡 Observer"
$wsReco.Cells.Item(27,7).Value = "➖ Neutre"

$wsReco.Cells.Item(28,1).Value = "BERNABE CI (BNBC)"
$wsReco.Cells.Item(28,2).Value = 2
$wsReco.Cells.Item(28,3).Value = 1
$wsReco.Cells.Item(28,4).Value = 4.69
$wsReco.Cells.Item(28,5).Value = 4.4
$wsReco.Cells.Item(28,6).Value = "🟡 Observer"
$wsReco.Cells.Item(28,7).Value = "👀 À surveiller"

$wsReco.Cells.Item(29,1).Value = "ORAGROUP TOGO (ORGT)"
$wsReco.Cells.Item(29,2).Value = 1
$wsReco.Cells.Item(29,3).Value = 0
$wsReco.Cells.Item(29,4).Value = 4.55
$wsReco.Cells.Item(29,5).Value = 4.55
$wsReco.Cells.Item(29,6).Value = "🟡 Observer"
$wsReco.Cells.Item(29,7).Value = "➖ Neutre"

$wsReco.Cells.Item(30,1).Value = "BICI CI (BICC)"
$wsReco.Cells.Item(30,2).Value = 1
$wsReco.Cells.Item(30,3).Value = 0
$wsReco.Cells.Item(30,4).Value = 3.97
$wsReco.Cells.Item(30,5).Value = 3.97
$wsReco.Cells.Item(30,6).Value = "🟡 Observer"
$wsReco.Cells.Item(30,7).Value = "➖ Neutre"

$wsReco.Cells.Item(31,1).Value = "SICOR CI (SICC)"
$wsReco.Cells.Item(31,2).Value = 1
$wsReco.Cells.Item(31,3).Value = 0
$wsReco.Cells.Item(31,4).Value = 2.94
$wsReco.Cells.Item(31,5).Value = 2.94
$wsReco.Cells.Item(31,6).Value = "🟡 Observer"
$wsReco.Cells.Item(31,7).Value = "➖ Neutre"

$wsReco.Cells.Item(32,1).Value = "TOTAL"
$wsReco.Cells.Item(32,2).Value = 0
$wsReco.Cells.Item(32,3).Value = 3
$wsReco.Cells.Item(32,4).Value = 0
$wsReco.Cells.Item(32,5).Value = 0
$wsReco.Cells.Item(32,6).Value = "🟡 Observer"
$wsReco.Cells.Item(32,7).Value = "➖ Neutre"

$wsReco.Cells.Item(33,1).Value = "UNIWAX CI (UNXC)"
$wsReco.Cells.Item(33,2).Value = 1
$wsReco.Cells.Item(33,3).Value = 1
$wsReco.Cells.Item(33,4).Value = -0.08
$wsReco.Cells.Item(33,5).Value = -5.08
$wsReco.Cells.Item(33,6).Value = "🟡 Observer"
$wsReco.Cells.Item(33,7).Value = "👀 À surveiller"

$wsReco.Cells.Item(34,1).Value = "TOTALENERGIES MARKETING SN (TTLS)"
$wsReco.Cells.Item(34,2).Value = 0
$wsReco.Cells.Item(34,3).Value = 1
$wsReco.Cells.Item(34,4).Value = -0.99
$wsReco.Cells.Item(34,5).Value = -0.99
$wsReco.Cells.Item(34,6).Value = "🟡 Observer"
$wsReco.Cells.Item(34,7).Value = "➖ Neutre"

$wsReco.Cells.Item(35,1).Value = "FILTISAC CI (FTSC)"
$wsReco.Cells.Item(35,2).Value = 0
$wsReco.Cells.Item(35,3).Value = 1
$wsReco.Cells.Item(35,4).Value = -1.2
$wsReco.Cells.Item(35,5).Value = -1.2
$wsReco.Cells.Item(35,6).Value = "🟡 Observer"
$wsReco.Cells.Item(35,7).Value = "➖ Neutre"

$wsReco.Cells.Item(36,1).Value = "AIR LIQUIDE CI (SIVC)"
$wsReco.Cells.Item(36,2).Value = 1
$wsReco.Cells.Item(36,3).Value = 1
$wsReco.Cells.Item(36,4).Value = -1.22
$wsReco.Cells.Item(36,5).Value = -5.63
$wsReco.Cells.Item(36,6).Value = "🟡 Observer"
$wsReco.Cells.Item(36,7).Value = "👀 À surveiller"

$wsReco.Cells.Item(37,1).Value = "BANK OF AFRICA ML (BOAM)"
$wsReco.Cells.Item(37,2).Value = 0
$wsReco.Cells.Item(37,3).Value = 1
$wsReco.Cells.Item(37,4).Value = -1.37
$wsReco.Cells.Item(37,5).Value = -1.37
$wsReco.Cells.Item(37,6).Value = "🟡 Observer"
$wsReco.Cells.Item(37,7).Value = "➖ Neutre"

$wsReco.Cells.Item(38,1).Value = "SMB CI (SMBC)"
$wsReco.Cells.Item(38,2).Value = 0
$wsReco.Cells.Item(38,3).Value = 1
$wsReco.Cells.Item(38,4).Value = -1.89
$wsReco.Cells.Item(38,5).Value = -1.89
$wsReco.Cells.Item(38,6).Value = "🟡 Observer"
$wsReco.Cells.Item(38,7).Value = "➖ Neutre"

$wsReco.Cells.Item(39,1).Value = "SITAB CI (STBC)"
$wsReco.Cells.Item(39,2).Value = 0
$wsReco.Cells.Item(39,3).Value = 1
$wsReco.Cells.Item(39,4).Value = -2.02
$wsReco.Cells.Item(39,5).Value = -2.02
$wsReco.Cells.Item(39,6).Value = "🟡 Observer"
$wsReco.Cells.Item(39,7).Value = "➖ Neutre"

$wsReco.Cells.Item(40,1).Value = "AFRICA GLOBAL LOGISTICS CI (SDSC)"
$wsReco.Cells.Item(40,2).Value = 0
$wsReco.Cells.Item(40,3).Value = 1
$wsReco.Cells.Item(40,4).Value = -2.68
$wsReco.Cells.Item(40,5).Value = -2.68
$wsReco.Cells.Item(40,6).Value = "🟡 Observer"
$wsReco.Cells.Item(40,7).Value = "➖ Neutre"

$wsReco.Cells.Item(41,1).Value = "ORANGE COTE D'IVOIRE (ORAC)"
$wsReco.Cells.Item(41,2).Value = 0
$wsReco.Cells.Item(41,3).Value = 1
$wsReco.Cells.Item(41,4).Value = -3.13
$wsReco.Cells.Item(41,5).Value = -3.13
$wsReco.Cells.Item(41,6).Value = "🟡 Observer"
$wsReco.Cells.Item(41,7).Value = "➖ Neutre"

$wsReco.Cells.Item(42,1).Value = "BANK OF AFRICA BF (BOABF)"
$wsReco.Cells.Item(42,2).Value = 0
$wsReco.Cells.Item(42,3).Value = 1
$wsReco.Cells.Item(42,4).Value = -3.31
$wsReco.Cells.Item(42,5).Value = -3.31
$wsReco.Cells.Item(42,6).Value = "🟡 Observer"
$wsReco.Cells.Item(42,7).Value = "➖ Neutre"

$wsReco.Cells.Item(43,1).Value = "ONATEL BF (ONTBF)"
$wsReco.Cells.Item(43,2).Value = 0
$wsReco.Cells.Item(43,3).Value = 1
$wsReco.Cells.Item(43,4).Value = -3.61
$wsReco.Cells.Item(43,5).Value = -3.61
$wsReco.Cells.Item(43,6).Value = "🟡 Observer"
$wsReco.Cells.Item(43,7).Value = "➖ Neutre"

$wsReco.Cells.Item(44,1).Value = "BANK OF AFRICA SENEGAL (BOAS)"
$wsReco.Cells.Item(44,2).Value = 0
$wsReco.Cells.Item(44,3).Value = 2
$wsReco.Cells.Item(44,4).Value = -4.24
$wsReco.Cells.Item(44,5).Value = -2.1
$wsReco.Cells.Item(44,6).Value = "🟡 Observer"
$wsReco.Cells.Item(44,7).Value = "➖ Neutre"

$wsReco.Cells.Item(45,1).Value = "VIVO ENERGY CI (SHEC)"
$wsReco.Cells.Item(45,2).Value = 0
$wsReco.Cells.Item(45,3).Value = 1
$wsReco.Cells.Item(45,4).Value = -4.8
$wsReco.Cells.Item(45,5).Value = -4.8
$wsReco.Cells.Item(45,6).Value = "🟡 Observer"
$wsReco.Cells.Item(45,7).Value = "➖ Neutre"

$wsReco.Cells.Item(46,1).Value = "LOTERIE NATIONALE DU BENIN (LNBB)"
$wsReco.Cells.Item(46,2).Value = 0
$wsReco.Cells.Item(46,3).Value = 2
$wsReco.Cells.Item(46,4).Value = -7.25
$wsReco.Cells.Item(46,5).Value = -3.5
$wsReco.Cells.Item(46,6).Value = "🟡 Observer"
$wsReco.Cells.Item(46,7).Value = "➖ Neutre"

$wsReco.Cells.Item(47,1).Value = "SICABLE CI (CABC)"
$wsReco.Cells.Item(47,2).Value = 0
$wsReco.Cells.Item(47,3).Value = 2
$wsReco.Cells.Item(47,4).Value = -11.68
$wsReco.Cells.Item(47,5).Value = -4.19
$wsReco.Cells.Item(47,6).Value = "🟡 Observer"
$wsReco.Cells.Item(47,7).Value = "➖ Neutre"

# --- Top_YTD sheet: rows 2-11 (B: Progression YTD (%)) ---
$wsYtd.Cells.Item(2,2).Value = 10560379.13
$wsYtd.Cells.Item(3,2).Value = 411766.78
$wsYtd.Cells.Item(4,2).Value = 379617.84
$wsYtd.Cells.Item(5,2).Value = 256849.95
$wsYtd.Cells.Item(6,2).Value = 170121.28
$wsYtd.Cells.Item(7,2).Value = 47791.58
$wsYtd.Cells.Item(8,2).Value = 43650.91
$wsYtd.Cells.Item(9,2).Value = 7904.76
$wsYtd.Cells.Item(10,2).Value = 3340.29
$wsYtd.Cells.Item(11,2).Value = 3327.36
